$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") "66.231.92"
$ws.Range("E2").Value = "  +4.62%  "

Set-TextCell $ws.Range("D3") "3.794.25"
$ws.Range("E3").Value = "  +7.32%  "

Set-TextCell $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.23%  "

Set-TextCell $ws.Range("D5") "429.26"
$ws.Range("E5").Value = "  +9.83%  "

Set-TextCell $ws.Range("D6") "139.67"
$ws.Range("E6").Value = "  +14.13%  "

$ws.Range("E7").Value = "  +5.68%  "

Set-TextCell $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +10.77%  "

$ws.Range("E10").Value = "  +3.49%  "

Set-TextCell $ws.Range("D11") "0.0000324"
$ws.Range("E11").Value = "  +0.32%  "

Set-TextCell $ws.Range("D12") "43.34"
$ws.Range("E12").Value = "  +12.57%  "

Set-TextCell $ws.Range("D13") "10.59"
$ws.Range("E13").Value = "  +16.37%  "

Set-TextCell $ws.Range("D14") "4.393.55"
$ws.Range("E14").Value = "  +7.42%  "

Set-TextCell $ws.Range("D15") "15.02"
$ws.Range("E15").Value = "  +14.37%  "

Set-TextCell $ws.Range("D17") "3.794.43"
$ws.Range("E17").Value = "  +7.69%  "

Set-TextCell $ws.Range("D18") "20.16"
$ws.Range("E18").Value = "  +8.44%  "

$ws.Range("E19").Value = "  +12.41%  "

Set-TextCell $ws.Range("D20") "66.385.27"
$ws.Range("E20").Value = "  +4.95%  "

Set-TextCell $ws.Range("D21") "410.51"
$ws.Range("E21").Value = "  +5.06%  "

$ws.Range("E22").Value = "  +10.79%  "

Set-TextCell $ws.Range("D23") "3.28"
$ws.Range("E23").Value = "  +13.92%  "

Set-TextCell $ws.Range("D24") "85.59"
$ws.Range("E24").Value = "  +5.65%  "

Set-TextCell $ws.Range("D25") "37.10"
$ws.Range("E25").Value = "  +10.77%  "

$ws.Range("E26").Value = "  +11.45%  "

Set-TextCell $ws.Range("D27") "9.67"
$ws.Range("E27").Value = "  +44.42%  "

Set-TextCell $ws.Range("D28") "9.86"
$ws.Range("E28").Value = "  +14.45%  "

Set-TextCell $ws.Range("D29") "5.42"
$ws.Range("E29").Value = "  -0.68%  "

Set-TextCell $ws.Range("D30") "13.88"
$ws.Range("E30").Value = "  +18.71%  "

Set-TextCell $ws.Range("D31") "699.82"
$ws.Range("E31").Value = "  +5.90%  "

$ws.Range("E32").Value = "  +18.06%  "

Set-TextCell $ws.Range("D33") "2.78"
$ws.Range("E33").Value = "  +7.50%  "

Set-TextCell $ws.Range("D34") "40.80"
$ws.Range("E34").Value = "  +12.22%  "

Set-TextCell $ws.Range("D36") "5.73"
$ws.Range("E36").Value = "  +41.30%  "

Set-TextCell $ws.Range("D37") "0.151"
$ws.Range("E37").Value = "  +2.39%  "

Set-TextCell $ws.Range("D38") "56.63"
$ws.Range("E38").Value = "  +6.36%  "

Set-TextCell $ws.Range("D39") "0.0478"
$ws.Range("E39").Value = "  +10.41%  "

$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell $ws.Range("D40") "2.88"
$ws.Range("E40").Value = "  +8.99%  "

Set-TextCell $ws.Range("D41") "0.0₃0680"
$ws.Range("E41").Value = "  +6.01%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D42") "0.141"
$ws.Range("E42").Value = "  +9.10%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D43") "2.58"
$ws.Range("E43").Value = "  +48.51%  "

$ws.Range("E44").Value = "  +0.56%  "

Set-TextCell $ws.Range("D45") "3.36"
$ws.Range("E45").Value = "  +10.24%  "

$ws.Range("E46").Value = "  +18.15%  "

$ws.Range("E47").Value = "  +5.74%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D48") "2.66"
$ws.Range("E48").Value = "  +8.67%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D49") "2.08"
$ws.Range("E49").Value = "  +6.70%  "

Set-TextCell $ws.Range("D50") "142.86"
$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("E51").Value = "  +7.78%  "
